$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Lower (E) and Upper (F) bound values reflecting recalculated
# uncertainty variation (pf calculation) for the TukeyHSD relative error
# confidence intervals.
$values = @(
    @{ Row = 2;  E = 0.01033416546889408; F = 0.06100621115665873 },
    @{ Row = 3;  E = 0.244137993531281;   F = 0.2947432627088423  },
    @{ Row = 4;  E = 0.07845460323788132; F = 0.1290047711323667  },
    @{ Row = 5;  E = 0.1851271578131059;  F = 0.2356773257075913  },
    @{ Row = 6;  E = 0.2101257750033669;  F = 0.2574151046112035  },
    @{ Row = 7;  E = 0.04444431890063329; F = 0.09167467884406195 },
    @{ Row = 8;  E = 0.1511168734758578;  F = 0.1983472334192865  },
    @{ Row = 9;  E = 0.1421315857977823;  F = 0.1892902960720929  },
    @{ Row = 10; E = 0.03545903122255778; F = 0.08261774149686829 },
    @{ Row = 11; E = 0.08312276603921226; F = 0.1302223431112369  }
)

foreach ($item in $values) {
    $ws.Cells.Item($item.Row, 5).Value = $item.E
    $ws.Cells.Item($item.Row, 6).Value = $item.F
}
